$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# employee_id, employee_name, department, absence_reason, absence_duration, absence_date, salary
$data = @(
    @(62526, "Vitor Pimenta", "Financeiro", "Doenca", 6, 45079, 2926.32),
    @(51123, "Heitor da Rocha", "Recursos Humanos", "Viagem de negocios", 2, 45084, 6183.82),
    @(11740, "Srta. Daniela Machado", "Financeiro", "Doenca", 5, 45099, 3031.34),
    @(13074, "Valentina Ribeiro", "Juridico", "Viagem de negocios", 8, 45097, 9136.42),
    @(72070, "Ana Luiza Mendonça", "Operacoes", "Problemas pessoais", 7, 45101, 8685.62),
    @(40263, "Maria Júlia Gomes", "Recursos Humanos", "Viagem de negocios", 7, 45099, 8438.73),
    @(34174, "Augusto Costa", "Financeiro", "Consulta medica", 2, 45080, 8124.38),
    @(67322, "Léo da Mata", "Atendimento ao Cliente", "Problemas pessoais", 6, 45105, 6351.58),
    @(49066, "Calebe Lopes", "Operacoes", "Doenca", 3, 45082, 6265.79),
    @(54342, "Arthur Miguel da Mota", "Vendas", "Viagem de negocios", 2, 45086, 7187.96)
)

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rec[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rec[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rec[5]
    $ws.Cells.Item($rowIndex, 7).Value = $rec[6]
    $rowIndex++
}
